# The workbook originally reported progress by generation number ("Gen")
# across 51 independent runs (Run 0 .. Run 50) plus a Mean column.
# It is being changed to report progress as a fraction of the maximum
# number of function evaluations ("MaxFES") across 50 runs (Run 0 .. Run 49)
# plus a recomputed Mean column (the Run 50 column is dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the first column header from "Gen" to "MaxFES" and replace
#    its values (generation counts) with the MaxFES fractions.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "MaxFES"

$maxfes = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $maxfes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $maxfes[$i]
}

# ---------------------------------------------------------------------
# 2) The old "Run 50" column (AZ) is being retired: turn it into the new
#    "Mean" column, recomputed over only the remaining 50 runs
#    (columns B..AY), then drop the old trailing "Mean" column (BA).
# ---------------------------------------------------------------------
$ws.Range("AZ1").Value = "Mean"

$newMean = @(
    62.75034427,
    59.38847983,
    31.80625809,
    1.04884747,
    0.31011575,
    0.26165116,
    0.24671892,
    0.2260322,
    0.21185301,
    0.191335,
    0.17864674,
    0.17219259,
    0.16734571
)
for ($i = 0; $i -lt $newMean.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 52).Value = $newMean[$i]
}

# Remove the now-redundant trailing "Mean" column (column BA / 53).
$ws.Columns.Item(53).Delete()
